$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: new columns AD (Wins), AE (Losses), AF (Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) so the
# new headers match the bold/centered/bordered look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the team record for every data row (2 through 54).
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 71  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 91  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
